$wb = $excel.ActiveWorkbook

# --- Sheet ALC: 30 cell updates ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 279.2  # ALC!H12: 399.33334 -> 279.2
$ws.Cells.Item(12, 10).Value = 348.75  # ALC!J12: 479 -> 348.75
$ws.Cells.Item(12, 12).Value = 348.75  # ALC!L12: 479 -> 348.75
$ws.Cells.Item(12, 14).Value = -688.75  # ALC!N12: -819 -> -688.75
$ws.Cells.Item(28, 8).Value = 484.3871  # ALC!H28: 516.069 -> 484.3871
$ws.Cells.Item(28, 9).Value = 306.26923  # ALC!I28: 319.0435 -> 306.26923
$ws.Cells.Item(28, 10).Value = 1410.6  # ALC!J28: 1271.3334 -> 1410.6
$ws.Cells.Item(28, 11).Value = 306.26923  # ALC!K28: 319.0435 -> 306.26923
$ws.Cells.Item(28, 12).Value = 1410.6  # ALC!L28: 1271.3334 -> 1410.6
$ws.Cells.Item(28, 13).Value = 178.73077  # ALC!M28: 165.9565 -> 178.73077
$ws.Cells.Item(28, 14).Value = -2380.6  # ALC!N28: -2241.3334 -> -2380.6
$ws.Cells.Item(33, 8).Value = 3788.4167  # ALC!H33: 3860.182 -> 3788.4167
$ws.Cells.Item(33, 9).Value = 3944.8333  # ALC!I33: 4134 -> 3944.8333
$ws.Cells.Item(33, 11).Value = 3944.8333  # ALC!K33: 4134 -> 3944.8333
$ws.Cells.Item(33, 13).Value = -3715.8333  # ALC!M33: -3905 -> -3715.8333
$ws.Cells.Item(51, 8).Value = 102499  # ALC!H51: 68000 -> 102499
$ws.Cells.Item(51, 9).Value = 4999  # ALC!I51: 2000.5 -> 4999
$ws.Cells.Item(51, 11).Value = 4999  # ALC!K51: 2000.5 -> 4999
$ws.Cells.Item(51, 13).Value = -4515  # ALC!M51: -1516.5 -> -4515
$ws.Cells.Item(108, 8).Value = 67500  # ALC!H108: 65776.664 -> 67500
$ws.Cells.Item(108, 10).Value = 60000  # ALC!J108: 57330 -> 60000
$ws.Cells.Item(108, 12).Value = 60000  # ALC!L108: 57330 -> 60000
$ws.Cells.Item(108, 14).Value = -67680  # ALC!N108: -65010 -> -67680
$ws.Cells.Item(137, 8).Value = 5511.3335  # ALC!H137: 5714.5806 -> 5511.3335
$ws.Cells.Item(137, 9).Value = 5910.143  # ALC!I137: 6094.55 -> 5910.143
$ws.Cells.Item(137, 10).Value = 4813.4165  # ALC!J137: 5023.727 -> 4813.4165
$ws.Cells.Item(137, 11).Value = 17730.429  # ALC!K137: 18283.65 -> 17730.429
$ws.Cells.Item(137, 12).Value = 14440.2495  # ALC!L137: 15071.181 -> 14440.2495
$ws.Cells.Item(137, 13).Value = -15180.429  # ALC!M137: -15733.65 -> -15180.429
$ws.Cells.Item(137, 14).Value = -19540.2495  # ALC!N137: -20171.181 -> -19540.2495

# --- Sheet ARM: 35 cell updates ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 16673831  # ARM!H32: 17864760 -> 16673831
$ws.Cells.Item(32, 9).Value = 25001846  # ARM!I32: 27779738 -> 25001846
$ws.Cells.Item(32, 11).Value = 25001846  # ARM!K32: 27779738 -> 25001846
$ws.Cells.Item(32, 13).Value = -25001559  # ARM!M32: -27779451 -> -25001559
$ws.Cells.Item(61, 8).Value = 83342440  # ARM!H61: 75008900 -> 83342440
$ws.Cells.Item(61, 10).Value = 62511250  # ARM!J61: 50010400 -> 62511250
$ws.Cells.Item(61, 12).Value = 62511250  # ARM!L61: 50010400 -> 62511250
$ws.Cells.Item(61, 14).Value = -62511674  # ARM!N61: -50010824 -> -62511674
$ws.Cells.Item(63, 8).Value = 4900.1665  # ARM!H63: 5135.5884 -> 4900.1665
$ws.Cells.Item(63, 9).Value = 1769.1  # ARM!I63: 1865.8889 -> 1769.1
$ws.Cells.Item(63, 11).Value = 1769.1  # ARM!K63: 1865.8889 -> 1769.1
$ws.Cells.Item(63, 13).Value = -1083.1  # ARM!M63: -1179.8889 -> -1083.1
$ws.Cells.Item(66, 8).Value = 4900.1665  # ARM!H66: 5135.5884 -> 4900.1665
$ws.Cells.Item(66, 9).Value = 1769.1  # ARM!I66: 1865.8889 -> 1769.1
$ws.Cells.Item(66, 11).Value = 8845.5  # ARM!K66: 9329.4445 -> 8845.5
$ws.Cells.Item(66, 13).Value = -5413.5  # ARM!M66: -5897.4445 -> -5413.5
$ws.Cells.Item(74, 8).Value = 20001128  # ARM!H74: 18572582 -> 20001128
$ws.Cells.Item(74, 10).Value = 3334502.2  # ARM!J74: 2501250.8 -> 3334502.2
$ws.Cells.Item(74, 12).Value = 3334502.2  # ARM!L74: 2501250.8 -> 3334502.2
$ws.Cells.Item(74, 14).Value = -3336250.2  # ARM!N74: -2502998.8 -> -3336250.2
$ws.Cells.Item(77, 8).Value = 20001128  # ARM!H77: 18572582 -> 20001128
$ws.Cells.Item(77, 10).Value = 3334502.2  # ARM!J77: 2501250.8 -> 3334502.2
$ws.Cells.Item(77, 12).Value = 16672511  # ARM!L77: 12506254 -> 16672511
$ws.Cells.Item(77, 14).Value = -16681247  # ARM!N77: -12514990 -> -16681247
$ws.Cells.Item(122, 8).Value = 2930.3447  # ARM!H122: 3025.5186 -> 2930.3447
$ws.Cells.Item(122, 9).Value = 2020.4736  # ARM!I122: 2038.5555 -> 2020.4736
$ws.Cells.Item(122, 10).Value = 4659.1  # ARM!J122: 4999.4443 -> 4659.1
$ws.Cells.Item(122, 11).Value = 6061.4208  # ARM!K122: 6115.666499999999 -> 6061.4208
$ws.Cells.Item(122, 12).Value = 13977.3  # ARM!L122: 14998.3329 -> 13977.3
$ws.Cells.Item(122, 13).Value = -3611.4208  # ARM!M122: -3665.666499999999 -> -3611.4208
$ws.Cells.Item(122, 14).Value = -18877.3  # ARM!N122: -19898.3329 -> -18877.3
$ws.Cells.Item(136, 8).Value = 83342440  # ARM!H136: 75008900 -> 83342440
$ws.Cells.Item(136, 10).Value = 62511250  # ARM!J136: 50010400 -> 62511250
$ws.Cells.Item(136, 12).Value = 187533750  # ARM!L136: 150031200 -> 187533750
$ws.Cells.Item(136, 14).Value = -187538850  # ARM!N136: -150036300 -> -187538850

# --- Sheet CRP: 8 cell updates ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 3090.2  # CRP!H134: 3203.7856 -> 3090.2
$ws.Cells.Item(134, 9).Value = 3090.2  # CRP!I134: 3203.7856 -> 3090.2
$ws.Cells.Item(134, 11).Value = 9270.599999999999  # CRP!K134: 9611.356800000001 -> 9270.599999999999
$ws.Cells.Item(134, 13).Value = -6735.599999999999  # CRP!M134: -7076.356800000001 -> -6735.599999999999
$ws.Cells.Item(140, 8).Value = 44100  # CRP!H140: 72050 -> 44100
$ws.Cells.Item(140, 10).Value = 0  # CRP!J140: 100000 -> 0
$ws.Cells.Item(140, 12).Value = 0  # CRP!L140: 100000 -> 0
$ws.Cells.Item(140, 14).Value = ""  # CRP!N140: -110360 -> (removed)

# --- Sheet CUL: 40 cell updates ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(86, 8).Value = 526  # CUL!H86: 525.875 -> 526
$ws.Cells.Item(86, 9).Value = 481.8  # CUL!I86: 481.6 -> 481.8
$ws.Cells.Item(86, 11).Value = 1445.4  # CUL!K86: 1444.8 -> 1445.4
$ws.Cells.Item(86, 13).Value = -259.4000000000001  # CUL!M86: -258.8000000000002 -> -259.4000000000001
$ws.Cells.Item(89, 8).Value = 526  # CUL!H89: 525.875 -> 526
$ws.Cells.Item(89, 9).Value = 481.8  # CUL!I89: 481.6 -> 481.8
$ws.Cells.Item(89, 11).Value = 4336.2  # CUL!K89: 4334.400000000001 -> 4336.2
$ws.Cells.Item(89, 13).Value = 1591.8  # CUL!M89: 1593.599999999999 -> 1591.8
$ws.Cells.Item(110, 8).Value = 11192.7  # CUL!H110: 12666.667 -> 11192.7
$ws.Cells.Item(110, 9).Value = 4013.5  # CUL!I110: 1000 -> 4013.5
$ws.Cells.Item(110, 10).Value = 12987.5  # CUL!J110: 15000 -> 12987.5
$ws.Cells.Item(110, 11).Value = 12040.5  # CUL!K110: 3000 -> 12040.5
$ws.Cells.Item(110, 12).Value = 38962.5  # CUL!L110: 45000 -> 38962.5
$ws.Cells.Item(110, 13).Value = -7950.5  # CUL!M110: 1090 -> -7950.5
$ws.Cells.Item(110, 14).Value = -47142.5  # CUL!N110: -53180 -> -47142.5
$ws.Cells.Item(113, 8).Value = 895.8333  # CUL!H113: 1067.375 -> 895.8333
$ws.Cells.Item(113, 9).Value = 679.8  # CUL!I113: 724.75 -> 679.8
$ws.Cells.Item(113, 10).Value = 978.9231  # CUL!J113: 1181.5834 -> 978.9231
$ws.Cells.Item(113, 11).Value = 2039.4  # CUL!K113: 2174.25 -> 2039.4
$ws.Cells.Item(113, 12).Value = 2936.7693  # CUL!L113: 3544.7502 -> 2936.7693
$ws.Cells.Item(113, 13).Value = 130.6000000000001  # CUL!M113: -4.25 -> 130.6000000000001
$ws.Cells.Item(113, 14).Value = -7276.7693  # CUL!N113: -7884.7502 -> -7276.7693
$ws.Cells.Item(117, 8).Value = 1735.6364  # CUL!H117: 1685.8182 -> 1735.6364
$ws.Cells.Item(117, 9).Value = 0  # CUL!I117: 775 -> 0
$ws.Cells.Item(117, 10).Value = 1735.6364  # CUL!J117: 1888.2222 -> 1735.6364
$ws.Cells.Item(117, 11).Value = 0  # CUL!K117: 2325 -> 0
$ws.Cells.Item(117, 12).Value = ""  # CUL!L117: 5664.6666 -> (removed)
$ws.Cells.Item(117, 13).Value = 5206.9092  # CUL!M117: 1117 -> 5206.9092
$ws.Cells.Item(117, 14).Value = -12090.9092  # CUL!N117: -12548.6666 -> -12090.9092
$ws.Cells.Item(131, 8).Value = 4598.2905  # CUL!H131: 4431.969 -> 4598.2905
$ws.Cells.Item(131, 10).Value = 4598.2905  # CUL!J131: 4431.969 -> 4598.2905
$ws.Cells.Item(131, 12).Value = 13794.8715  # CUL!L131: 13295.907 -> 13794.8715
$ws.Cells.Item(131, 14).Value = -23874.8715  # CUL!N131: -23375.907 -> -23874.8715
$ws.Cells.Item(132, 8).Value = 1956.7  # CUL!H132: 2172.125 -> 1956.7
$ws.Cells.Item(132, 9).Value = 1235  # CUL!I132: 1480 -> 1235
$ws.Cells.Item(132, 10).Value = 2137.125  # CUL!J132: 2271 -> 2137.125
$ws.Cells.Item(132, 11).Value = 11115  # CUL!K132: 13320 -> 11115
$ws.Cells.Item(132, 12).Value = 19234.125  # CUL!L132: 20439 -> 19234.125
$ws.Cells.Item(132, 13).Value = -8585  # CUL!M132: -10790 -> -8585
$ws.Cells.Item(132, 14).Value = -24294.125  # CUL!N132: -25499 -> -24294.125

# --- Sheet GSM: 15 cell updates ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 239.09091  # GSM!H2: 282.22223 -> 239.09091
$ws.Cells.Item(2, 9).Value = 62.5  # GSM!I2: 56.666668 -> 62.5
$ws.Cells.Item(2, 10).Value = 340  # GSM!J2: 395 -> 340
$ws.Cells.Item(2, 11).Value = 62.5  # GSM!K2: 56.666668 -> 62.5
$ws.Cells.Item(2, 12).Value = 340  # GSM!L2: 395 -> 340
$ws.Cells.Item(2, 13).Value = 50.5  # GSM!M2: 56.333332 -> 50.5
$ws.Cells.Item(2, 14).Value = -566  # GSM!N2: -621 -> -566
$ws.Cells.Item(102, 8).Value = 2901.2  # GSM!H102: 3497.4736 -> 2901.2
$ws.Cells.Item(102, 9).Value = 2193.375  # GSM!I102: 2901.6 -> 2193.375
$ws.Cells.Item(102, 11).Value = 2193.375  # GSM!K102: 2901.6 -> 2193.375
$ws.Cells.Item(102, 13).Value = -571.375  # GSM!M102: -1279.6 -> -571.375
$ws.Cells.Item(132, 8).Value = 34490380  # GSM!H132: 35722104 -> 34490380
$ws.Cells.Item(132, 9).Value = 47622584  # GSM!I132: 50003610 -> 47622584
$ws.Cells.Item(132, 11).Value = 142867752  # GSM!K132: 150010830 -> 142867752
$ws.Cells.Item(132, 13).Value = -142865222  # GSM!M132: -150008300 -> -142865222

# --- Sheet LTW: 25 cell updates ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 3331.9285  # LTW!H46: 3242.077 -> 3331.9285
$ws.Cells.Item(46, 10).Value = 4881.2  # LTW!J46: 4976.5 -> 4881.2
$ws.Cells.Item(46, 12).Value = 4881.2  # LTW!L46: 4976.5 -> 4881.2
$ws.Cells.Item(46, 14).Value = -5257.2  # LTW!N46: -5352.5 -> -5257.2
$ws.Cells.Item(60, 8).Value = 29833.334  # LTW!H60: 36000 -> 29833.334
$ws.Cells.Item(60, 9).Value = 29500  # LTW!I60: 36000 -> 29500
$ws.Cells.Item(60, 10).Value = 30000  # LTW!J60: 0 -> 30000
$ws.Cells.Item(60, 11).Value = 29500  # LTW!K60: 36000 -> 29500
$ws.Cells.Item(60, 12).Value = 30000  # LTW!L60: 0 -> 30000
$ws.Cells.Item(60, 13).Value = -28991  # LTW!M60: -35491 -> -28991
$ws.Cells.Item(60, 14).Value = -31018  # LTW!N60: (new) -> -31018
$ws.Cells.Item(93, 8).Value = 58824880  # LTW!H93: 71430120 -> 58824880
$ws.Cells.Item(93, 9).Value = 125000860  # LTW!I93: 142858130 -> 125000860
$ws.Cells.Item(93, 10).Value = 1782.8889  # LTW!J93: 2124.5715 -> 1782.8889
$ws.Cells.Item(93, 11).Value = 125000860  # LTW!K93: 142858130 -> 125000860
$ws.Cells.Item(93, 12).Value = 1782.8889  # LTW!L93: 2124.5715 -> 1782.8889
$ws.Cells.Item(93, 13).Value = -124999612  # LTW!M93: -142856882 -> -124999612
$ws.Cells.Item(93, 14).Value = -4278.8889  # LTW!N93: -4620.5715 -> -4278.8889
$ws.Cells.Item(136, 8).Value = 74289.63  # LTW!H136: 71277.13 -> 74289.63
$ws.Cells.Item(136, 9).Value = 11105.23  # LTW!I136: 11743.833 -> 11105.23
$ws.Cells.Item(136, 10).Value = 165556  # LTW!J136: 136222.55 -> 165556
$ws.Cells.Item(136, 11).Value = 33315.69  # LTW!K136: 35231.499 -> 33315.69
$ws.Cells.Item(136, 12).Value = 496668  # LTW!L136: 408667.65 -> 496668
$ws.Cells.Item(136, 13).Value = -30765.69  # LTW!M136: -32681.499 -> -30765.69
$ws.Cells.Item(136, 14).Value = -501768  # LTW!N136: -413767.65 -> -501768

# --- Sheet WVR: 20 cell updates ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(69, 8).Value = 29757  # WVR!H69: 29635.5 -> 29757
$ws.Cells.Item(69, 10).Value = 29757  # WVR!J69: 29635.5 -> 29757
$ws.Cells.Item(69, 12).Value = 29757  # WVR!L69: 29635.5 -> 29757
$ws.Cells.Item(69, 14).Value = -31255  # WVR!N69: -31133.5 -> -31255
$ws.Cells.Item(72, 8).Value = 29757  # WVR!H72: 29635.5 -> 29757
$ws.Cells.Item(72, 10).Value = 29757  # WVR!J72: 29635.5 -> 29757
$ws.Cells.Item(72, 12).Value = 89271  # WVR!L72: 88906.5 -> 89271
$ws.Cells.Item(72, 14).Value = -96759  # WVR!N72: -96394.5 -> -96759
$ws.Cells.Item(81, 8).Value = 6229.769  # WVR!H81: 5505.067 -> 6229.769
$ws.Cells.Item(81, 9).Value = 1713.1428  # WVR!I81: 1509 -> 1713.1428
$ws.Cells.Item(81, 11).Value = 3426.2856  # WVR!K81: 3018 -> 3426.2856
$ws.Cells.Item(81, 13).Value = -2365.2856  # WVR!M81: -1957 -> -2365.2856
$ws.Cells.Item(84, 8).Value = 6229.769  # WVR!H84: 5505.067 -> 6229.769
$ws.Cells.Item(84, 9).Value = 1713.1428  # WVR!I84: 1509 -> 1713.1428
$ws.Cells.Item(84, 11).Value = 17131.428  # WVR!K84: 15090 -> 17131.428
$ws.Cells.Item(84, 13).Value = -11827.428  # WVR!M84: -9786 -> -11827.428
$ws.Cells.Item(136, 8).Value = 6473.8  # WVR!H136: 6774.1333 -> 6473.8
$ws.Cells.Item(136, 10).Value = 4475  # WVR!J136: 5038.125 -> 4475
$ws.Cells.Item(136, 12).Value = 13425  # WVR!L136: 15114.375 -> 13425
$ws.Cells.Item(136, 14).Value = -18525  # WVR!N136: -20214.375 -> -18525
